$d = $word.ActiveDocument

# --- Edit 1: "Administrador, Contratante" -> "Administrativo, Contratante" ---
$d.Content.Find.Execute(
    "3 tipos de usuarios: Administrador, Contratante y ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "3 tipos de usuarios: Administrativo, Contratante y ", 2) | Out-Null

# --- Edit 2: ". Si es administrador" -> ". Si es administrativo" ---
$d.Content.Find.Execute(
    ". Si es administrador",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ". Si es administrativo", 2) | Out-Null

# --- Edit 3: remove the stray _GoBack bookmark that currently splits
#     "parta"/"mento" inside "Departamento" (it gets relocated in Edit 6) ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Edit 4: "Un administrador puede aceptar" -> "Un usuario administrativo puede aceptar" ---
$d.Content.Find.Execute(
    "Un administrador puede aceptar",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Un usuario administrativo puede aceptar", 2) | Out-Null

# --- Edit 5: "realizado por el administrador, este" -> "realizado por el ente administrativo, este" ---
$d.Content.Find.Execute(
    "realizado por el administrador, este",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "realizado por el ente administrativo, este", 2) | Out-Null

# --- Edit 6: "El administrador deberá dejar" -> "El usuario administrativo deberá dejar" ---
$d.Content.Find.Execute(
    "El administrador deber",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El usuario administrativo deber", 2) | Out-Null

# Relocate the _GoBack bookmark to sit right before "deberá" in that sentence,
# matching the edited document (bookmark now falls between
# "usuario administrativo " and "deberá dejar algunos Estudiantes...").
$rng = $d.Content
$rng.Find.Execute("deberá dejar algunos Estudiantes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target = $d.Range($rng.Start, $rng.Start)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
